$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1671
$ws1.Range("F9").Value = 3507
$ws1.Range("F14").Value = 884
$ws1.Range("F16").Value = 1261
$ws1.Range("F17").Value = 1782
$ws1.Range("F24").Value = 4247

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 119

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F13").Value = 792
$ws3.Range("F14").Value = 194

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 792
$ws4.Range("F20").Value = 884
$ws4.Range("F23").Value = 1261
$ws4.Range("F30").Value = 1782
$ws4.Range("F35").Value = 119
$ws4.Range("F36").Value = 119
$ws4.Range("F44").Value = 4247
